# ============================================================
# Ecuador LigaPro Serie A workbook update
# - Adds new team "Imbabura"
# - Swaps mis-ordered match rows (135/136, 139/140, 143/144)
# - Updates odds for an existing upcoming fixture
# - Inserts one newly completed match and five new upcoming fixtures
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$allCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Set-RowValuesByCols($ws, $rowNum, $colLetters, $values) {
    for ($i = 0; $i -lt $colLetters.Count; $i++) {
        $cellRef = "$($colLetters[$i])$rowNum"
        $ws.Range($cellRef).Value = $values[$i]
    }
}

# --- Fix rows 135/136, 139/140, 143/144: content was swapped between the two rows ---
$v135 = @(133, 7483188, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45256.83333333334, "Gualaceo SC", "Emelec", 0, 2, "A", 3.6, 3.3, 2.05, 2.6, 3.25, 2.75, 0, 1.8, 2, 2.5, 1.975, 1.825, -1, -1, 1.75, -1, 1, -1, 0.825)
Set-RowValuesByCols $ws 135 $allCols $v135
$v136 = @(134, 7483306, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45256.83333333334, "Tecnico Universitario", "Club Atletico Libertad", 1, 1, "D", 1.5, 4.333, 5.75, 1.533, 4.2, 5.5, -1, 1.925, 1.875, 2.25, 1.8, 2, -1, 3.2, -1, -1, 0.875, -0.5, 0.5)
Set-RowValuesByCols $ws 136 $allCols $v136
$v139 = @(137, 7528849, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45262.70833333334, "Guayaquil City", "Gualaceo SC", 0, 2, "A", 1.833, 3.5, 3.75, 2.15, 3.4, 3, -0.25, 1.825, 1.975, 2.5, 1.85, 1.95, -1, -1, 2, -1, 0.9750000000000001, -1, 0.95)
Set-RowValuesByCols $ws 139 $allCols $v139
$v140 = @(138, 7528859, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45262.70833333334, "Club Atletico Libertad", "Cumbaya FC", 3, 1, "H", 1.727, 3.5, 4.333, 1.4, 4.2, 7, -1.25, 2, 1.8, 2.5, 1.95, 1.85, 0.3999999999999999, -1, -1, 1, -1, 0.95, -1)
Set-RowValuesByCols $ws 140 $allCols $v140
$v143 = @(141, 7528857, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45263.83333333334, "Universidad Catolica del Ecuador", "Barcelona Guayaquil", 0, 1, "A", 1.533, 4, 5.5, 1.5, 4.333, 5.25, -1, 1.8, 2, 3, 1.975, 1.825, -1, -1, 4.25, -1, 1, -1, 0.825)
Set-RowValuesByCols $ws 143 $allCols $v143
$v144 = @(142, 7528852, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45263.83333333334, "Delfin SC", "Tecnico Universitario", 2, 2, "D", 2.1, 3.4, 3.1, 2.1, 3.4, 3.1, -0.25, 1.8, 2, 2.25, 1.9, 1.9, -1, 2.4, -1, -0.5, 0.5, 0.8999999999999999, -1)
Set-RowValuesByCols $ws 144 $allCols $v144
# --- A new completed match is inserted as row 148 (pushes old row148 down to 149) ---
$ws.Rows.Item(148).Insert()

$partialCols = @("A","B","C","D","E","F","G","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

$v148 = @(146, 7773060, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45352.875, "Delfin SC", "Independiente del Valle", 0, 1, "A", 3.1, 3.1, 2.3, 3.75, 3.2, 2, 0.5, 1.825, 1.975, 2.25, 2, 1.8, -1, -1, 1, -1, 0.9750000000000001, -1, 0.8)
Set-RowValuesByCols $ws 148 $allCols $v148
# copy number/border formatting for the id (A) and date (E) cells from the row above
$ws.Range("A147").Copy()
$ws.Range("A148").PasteSpecial(-4122)
$ws.Range("E147").Copy()
$ws.Range("E148").PasteSpecial(-4122)

# --- Row 149 (the shifted former row 148) only had its closing odds updated ---
$ws.Range("A149").Value = 147
$ws.Range("N149").Value = 1.833
$ws.Range("O149").Value = 3.5
$ws.Range("R149").Value = 1.825
$ws.Range("S149").Value = 1.975
$ws.Range("U149").Value = 1.95
$ws.Range("V149").Value = 1.85

# --- Five new upcoming fixtures appended as rows 150-154 (no result yet) ---
$v150 = @(148, 7773461, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45353.72916666666, "SD Aucas", "Orense", 1.7, 3.75, 4.75, 1.5, 4.2, 6, -0.75, 1.85, 1.95, 2.75, 1.975, 1.825, 0, 0, 0, 0, 0)
Set-RowValuesByCols $ws 150 $partialCols $v150
$ws.Range("A149").Copy()
$ws.Range("A150").PasteSpecial(-4122)
$ws.Range("E149").Copy()
$ws.Range("E150").PasteSpecial(-4122)
$v151 = @(149, 7773458, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45353.83333333334, "Deportivo Cuenca", "Emelec", 2.375, 3.25, 2.9, 2.7, 3.2, 2.55, 0, 1.95, 1.85, 2.25, 1.9, 1.9, 0, 0, 0, 0, 0)
Set-RowValuesByCols $ws 151 $partialCols $v151
$ws.Range("A150").Copy()
$ws.Range("A151").PasteSpecial(-4122)
$ws.Range("E150").Copy()
$ws.Range("E151").PasteSpecial(-4122)
$v152 = @(150, 7798121, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45354.625, "Universidad Catolica del Ecuador", "El Nacional", 1.75, 3.5, 4.333, 1.615, 3.75, 5.25, -0.75, 1.775, 2.025, 2.75, 1.85, 1.95, 0, 0, 0, 0, 0)
Set-RowValuesByCols $ws 152 $partialCols $v152
$ws.Range("A151").Copy()
$ws.Range("A152").PasteSpecial(-4122)
$ws.Range("E151").Copy()
$ws.Range("E152").PasteSpecial(-4122)
$v153 = @(151, 7773460, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45354.72916666666, "Tecnico Universitario", "Cumbaya FC", 1.533, 4, 6, 1.45, 4.2, 6.5, -1, 1.85, 1.95, 2.25, 1.9, 1.9, 0, 0, 0, 0, 0)
Set-RowValuesByCols $ws 153 $partialCols $v153
$ws.Range("A152").Copy()
$ws.Range("A153").PasteSpecial(-4122)
$ws.Range("E152").Copy()
$ws.Range("E153").PasteSpecial(-4122)
$v154 = @(152, 7773780, "Ecuador LigaPro Serie A", "Ecuador LigaPro Serie A", 45354.83333333334, "Barcelona Guayaquil", "Imbabura", 1.2, 6, 11, 1.222, 6, 10, -1.75, 1.95, 1.85, 2.75, 1.9, 1.9, 0, 0, 0, 0, 0)
Set-RowValuesByCols $ws 154 $partialCols $v154
$ws.Range("A153").Copy()
$ws.Range("A154").PasteSpecial(-4122)
$ws.Range("E153").Copy()
$ws.Range("E154").PasteSpecial(-4122)
Write-Host "Edit complete."
